# iron_native COM-interop edit script
# ------------------------------------------------------------------
# Source data for this sheet used to report eval_avg_ep_rew /
# eval_max_ep_rew / eval_min_ep_rew as three separate leading columns.
# The run now reports a single "total_time_taken(m)" column instead,
# a new "O_num_agents" hyper-param column was added, and four rows of
# actual hyper-parameter-tuning run results were appended.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row changes -------------------------------------------
# Collapse the old A1:C1 "eval_avg_ep_rew / eval_max_ep_rew /
# eval_min_ep_rew" columns into a single "total_time_taken(m)" column:
# delete B and C, then rename A1.
$ws.Range("B1:C1").EntireColumn.Delete()
$ws.Range("A1").Value = "total_time_taken(m)"

# Insert a new "O_num_agents" column right before "O_reach_goal_rew"
# (after the shift above, O_reach_goal_rew now sits at AK1).
$ws.Range("AK1").EntireColumn.Insert()
$ws.Range("AK1").Value = "O_num_agents"

# --- New data rows ---------------------------------------------------
# Columns, in sheet order (R = M_edge_feat_size and T = M_edge_hidden_size
# are intentionally left blank in every row, matching the source data).
$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q", `
             "S","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH", `
             "AI","AJ","AK","AL","AM")

$row2 = @(0.6579903721809387, 0.303030303030303, 3.41, 9.24, 20, 128, 10, 100, 10, $false, 1, 4, 6, $true, 40, 202, 42, `
          32, 32, 1, 2, $true, 32, 0.95, 0.01, 32, 10, 0.2, 0.9, 0.001, $true, 5, `
          0.8, 0.9, 8, 10, -5)

$row3 = @(0.603220264116923, 0.2626262626262627, 3.37, 9.15, 20, 128, 10, 100, 10, $false, 1, 4, 6, $true, 40, 202, 42, `
          32, 32, 1, 2, $true, 32, 0.95, 0.01, 32, 10, 0.2, 0.9, 0.001, $true, 5, `
          0.8, 0.9, 8, 10, -5)

$row4 = @(0.5995295246442159, 0.3636363636363636, 3.35, 8.98, 20, 128, 10, 100, 10, $false, 1, 4, 6, $true, 40, 202, 42, `
          32, 32, 1, 2, $true, 32, 0.95, 0.01, 32, 10, 0.2, 0.95, 0.001, $true, 5, `
          0.2, 0.75, 8, 10, -5)

$row5 = @(0.5993634541829427, 0.4444444444444444, 3.48, 8.82, 20, 128, 10, 100, 10, $false, 1, 4, 6, $true, 40, 202, 42, `
          32, 32, 1, 2, $true, 32, 0.95, 0.01, 32, 10, 0.2, 0.95, 0.001, $true, 5, `
          0.2, 0.75, 8, 10, -5)

$dataRows = @($row2, $row3, $row4, $row5)

$rowIndex = 2
foreach ($rowValues in $dataRows) {
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $addr = "$($columns[$i])$rowIndex"
        $ws.Range($addr).Value = $rowValues[$i]
    }
    $rowIndex = $rowIndex + 1
}
